$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1544.3334
$ws.Range("I11").Value = 1544.3334
$ws.Range("K11").Value = 1544.3334
$ws.Range("M11").Value = -1404.3334
# Row 41
$ws.Range("H41").Value = 7812891.5
$ws.Range("I41").Value = 13889112
$ws.Range("K41").Value = 13889112
$ws.Range("M41").Value = -13888672
# Row 53
$ws.Range("H53").Value = 3316.05
$ws.Range("I53").Value = 2990.4546
$ws.Range("J53").Value = 3714
$ws.Range("K53").Value = 2990.4546
$ws.Range("L53").Value = 3714
$ws.Range("M53").Value = -2353.4546
$ws.Range("N53").Value = -4988
# Row 62
$ws.Range("H62").Value = 88880
$ws.Range("J62").Value = 88880
$ws.Range("L62").Value = 88880
$ws.Range("N62").Value = -90128
# Row 65
$ws.Range("H65").Value = 88880
$ws.Range("J65").Value = 88880
$ws.Range("L65").Value = 444400
$ws.Range("N65").Value = -450640
# Row 86
$ws.Range("H86").Value = 62138660
$ws.Range("I86").Value = 93753064
$ws.Range("K86").Value = 93753064
$ws.Range("M86").Value = -93751941
# Row 89
$ws.Range("H89").Value = 62138660
$ws.Range("I89").Value = 93753064
$ws.Range("K89").Value = 468765320
$ws.Range("M89").Value = -468759704
# Row 98
$ws.Range("H98").Value = 4035.5
$ws.Range("I98").Value = 5336.4287
$ws.Range("K98").Value = 5336.4287
$ws.Range("M98").Value = -3838.4287
# Row 103
$ws.Range("H103").Value = 1149.3125
$ws.Range("J103").Value = 1296.9231
$ws.Range("L103").Value = 3890.7693
$ws.Range("N103").Value = -5062.7693
# Row 122
$ws.Range("H122").Value = 4035.5
$ws.Range("I122").Value = 5336.4287
$ws.Range("K122").Value = 16009.2861
$ws.Range("M122").Value = -13559.2861
# Row 132
$ws.Range("H132").Value = 1808.4445
$ws.Range("I132").Value = 1872.64
$ws.Range("K132").Value = 5617.92
$ws.Range("M132").Value = -3087.92
# Row 138
$ws.Range("H138").Value = 2623.1304
$ws.Range("J138").Value = 2654.3076
$ws.Range("L138").Value = 7962.9228
$ws.Range("N138").Value = -18242.9228

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1789671.4
$ws.Range("I32").Value = 1789671.4
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1789671.4
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1789384.4
$ws.Range("N32").ClearContents()
# Row 61
$ws.Range("H61").Value = 2995.9019
$ws.Range("I61").Value = 1936.7142
$ws.Range("K61").Value = 1936.7142
$ws.Range("M61").Value = -1724.7142
# Row 88
$ws.Range("H88").Value = 1514.3889
$ws.Range("I88").Value = 1342.4
$ws.Range("J88").Value = 1580.5385
$ws.Range("K88").Value = 1342.4
$ws.Range("L88").Value = 1580.5385
$ws.Range("M88").Value = -936.4000000000001
$ws.Range("N88").Value = -2392.5385
# Row 91
$ws.Range("H91").Value = 1514.3889
$ws.Range("I91").Value = 1342.4
$ws.Range("J91").Value = 1580.5385
$ws.Range("K91").Value = 1342.4
$ws.Range("L91").Value = 1580.5385
$ws.Range("M91").Value = 61.59999999999991
$ws.Range("N91").Value = -4388.538500000001
# Row 132
$ws.Range("H132").Value = 4862.878
$ws.Range("I132").Value = 4397.407
$ws.Range("K132").Value = 13192.221
$ws.Range("M132").Value = -10662.221
# Row 136
$ws.Range("H136").Value = 2995.9019
$ws.Range("I136").Value = 1936.7142
$ws.Range("K136").Value = 5810.142599999999
$ws.Range("M136").Value = -3260.142599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 39840.355
$ws.Range("I86").Value = 48705.953
$ws.Range("J86").Value = 7333.1665
$ws.Range("K86").Value = 48705.953
$ws.Range("L86").Value = 7333.1665
$ws.Range("M86").Value = -47582.953
$ws.Range("N86").Value = -9579.166499999999
# Row 89
$ws.Range("H89").Value = 39840.355
$ws.Range("I89").Value = 48705.953
$ws.Range("J89").Value = 7333.1665
$ws.Range("K89").Value = 243529.765
$ws.Range("L89").Value = 36665.8325
$ws.Range("M89").Value = -237913.765
$ws.Range("N89").Value = -47897.8325
# Row 134
$ws.Range("H134").Value = 4323.8716
$ws.Range("I134").Value = 1709.7
$ws.Range("J134").Value = 7075.6313
$ws.Range("K134").Value = 5129.1
$ws.Range("L134").Value = 21226.8939
$ws.Range("M134").Value = -2594.1
$ws.Range("N134").Value = -26296.8939

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8340748.5
$ws.Range("I31").Value = 2506.2727
$ws.Range("K31").Value = 2506.2727
$ws.Range("M31").Value = -2211.2727
# Row 34
$ws.Range("H34").Value = 8340748.5
$ws.Range("I34").Value = 2506.2727
$ws.Range("K34").Value = 2506.2727
$ws.Range("M34").Value = -2304.2727
# Row 58
$ws.Range("H58").Value = 7406.4194
$ws.Range("I58").Value = 2022.2727
$ws.Range("J58").Value = 10367.7
$ws.Range("K58").Value = 2022.2727
$ws.Range("L58").Value = 10367.7
$ws.Range("M58").Value = -1819.2727
$ws.Range("N58").Value = -10773.7
# Row 136
$ws.Range("H136").Value = 7406.4194
$ws.Range("I136").Value = 2022.2727
$ws.Range("J136").Value = 10367.7
$ws.Range("K136").Value = 6066.8181
$ws.Range("L136").Value = 31103.1
$ws.Range("M136").Value = -3516.8181
$ws.Range("N136").Value = -36203.10000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 402002020
$ws.Range("I4").Value = 336666720
$ws.Range("K4").Value = 1010000160
$ws.Range("M4").Value = -1010000048
# Row 50
$ws.Range("H50").Value = 83333530
$ws.Range("I50").Value = 83333530
$ws.Range("K50").Value = 250000590
$ws.Range("M50").Value = -250000109
# Row 53
$ws.Range("H53").Value = 83333530
$ws.Range("I53").Value = 83333530
$ws.Range("K53").Value = 250000590
$ws.Range("M53").Value = -250000109
# Row 80
$ws.Range("H80").Value = 52636100
$ws.Range("J80").Value = 71434130
$ws.Range("L80").Value = 214302390
$ws.Range("N80").Value = -214304262
# Row 83
$ws.Range("H83").Value = 52636100
$ws.Range("J83").Value = 71434130
$ws.Range("L83").Value = 642907170
$ws.Range("N83").Value = -642916530

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 33333.332
$ws.Range("J47").Value = 33333.332
$ws.Range("L47").Value = 33333.332
$ws.Range("N47").Value = -34469.332
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 80
$ws.Range("H80").Value = 3979
$ws.Range("I80").Value = 3335.3333
$ws.Range("J80").Value = 4365.2
$ws.Range("K80").Value = 3335.3333
$ws.Range("L80").Value = 4365.2
$ws.Range("M80").Value = -2337.3333
$ws.Range("N80").Value = -6361.2
# Row 83
$ws.Range("H83").Value = 3979
$ws.Range("I83").Value = 3335.3333
$ws.Range("J83").Value = 4365.2
$ws.Range("K83").Value = 16676.6665
$ws.Range("L83").Value = 21826
$ws.Range("M83").Value = -11684.6665
$ws.Range("N83").Value = -31810
# Row 122
$ws.Range("H122").Value = 60075.332
$ws.Range("I122").Value = 94604.73
$ws.Range("K122").Value = 283814.19
$ws.Range("M122").Value = -281364.19
# Row 126
$ws.Range("H126").Value = 2628.5
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3014
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9042
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -13982

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 29
$ws.Range("H29").Value = 1183
$ws.Range("I29").Value = 1183
$ws.Range("K29").Value = 1183
$ws.Range("M29").Value = -888
# Row 40
$ws.Range("H40").Value = 6724.9062
$ws.Range("I40").Value = 5019.2
$ws.Range("K40").Value = 5019.2
$ws.Range("M40").Value = -4883.2
# Row 61
$ws.Range("H61").Value = 8109.933
$ws.Range("I61").Value = 6742.3335
$ws.Range("K61").Value = 6742.3335
$ws.Range("M61").Value = -6540.3335
# Row 113
$ws.Range("H113").Value = 8109.933
$ws.Range("I113").Value = 6742.3335
$ws.Range("K113").Value = 6742.3335
$ws.Range("M113").Value = -4572.3335
# Row 132
$ws.Range("H132").Value = 4972.7446
$ws.Range("I132").Value = 2989.2693
$ws.Range("K132").Value = 8967.8079
$ws.Range("M132").Value = -6437.8079

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 16209.875
$ws.Range("I54").Value = 14946.5
$ws.Range("K54").Value = 14946.5
$ws.Range("M54").Value = -14426.5
# Row 122
$ws.Range("H122").Value = 146927.53
$ws.Range("I122").Value = 212549.27
$ws.Range("J122").Value = 8392.777
$ws.Range("K122").Value = 637647.8099999999
$ws.Range("L122").Value = 25178.331
$ws.Range("M122").Value = -635197.8099999999
$ws.Range("N122").Value = -30078.331
# Row 132
$ws.Range("H132").Value = 5446.4688
$ws.Range("I132").Value = 6658.9
$ws.Range("J132").Value = 3425.75
$ws.Range("K132").Value = 19976.7
$ws.Range("L132").Value = 10277.25
$ws.Range("M132").Value = -17446.7
$ws.Range("N132").Value = -15337.25
# Row 136
$ws.Range("H136").Value = 52588.59
$ws.Range("I136").Value = 2242.75
$ws.Range("K136").Value = 6728.25
$ws.Range("M136").Value = -4178.25
